# Adds OBJs 106, 107, 108 and their associated source files to the
# Eve web catalog workbook (Objects / Images / Albums / Types and Tags).

$wb = $excel.ActiveWorkbook

$wsObjects = $wb.Worksheets.Item("Objects")
$wsImages  = $wb.Worksheets.Item("Images")
$wsAlbums  = $wb.Worksheets.Item("Albums")
$wsTypes   = $wb.Worksheets.Item("Types and Tags")

# --- Objects sheet: fill in rows 45-47 (CAT NO. already present in col A) ---
$wsObjects.Range("B45").Value = "2015-08-09-obj000106-01.textile"
$wsObjects.Range("C45").Value = "art"
$wsObjects.Range("D45").Value = "artworks"
$wsObjects.Range("E45").Value = "Fabric/hanging"
$wsObjects.Range("F45").Value = "Untitled"
$wsObjects.Range("G45").Value = "pic000209"
$wsObjects.Range("H45").Value = "alb000106"

$wsObjects.Range("B46").Value = "2015-08-09-obj000107-01.textile"
$wsObjects.Range("C46").Value = "art"
$wsObjects.Range("D46").Value = "artworks"
$wsObjects.Range("E46").Value = "Fabric/hanging"
$wsObjects.Range("F46").Value = "Untitled"
$wsObjects.Range("G46").Value = "pic000210"
$wsObjects.Range("H46").Value = "alb000107"

$wsObjects.Range("B47").Value = "2015-08-09-obj000108-01.textile"
$wsObjects.Range("C47").Value = "art"
$wsObjects.Range("D47").Value = "artworks"
$wsObjects.Range("E47").Value = "Fabric/hanging"
$wsObjects.Range("F47").Value = "Untitled"
$wsObjects.Range("G47").Value = "pic000211"
$wsObjects.Range("H47").Value = "alb000108"

# --- Images sheet: new rows for the three new pics ---
$wsImages.Range("A150").Value = "pic000209"
$wsImages.Range("B150").Value = "2015-08-09-pic000209.textile"
$wsImages.Range("C150").Value = "S-3-0031.jpg"
$wsImages.Range("I150").Value = "alb000106"

$wsImages.Range("A151").Value = "pic000210"
$wsImages.Range("B151").Value = "2015-08-09-pic000210.textile"
$wsImages.Range("C151").Value = "S-14-0021.jpg"
$wsImages.Range("I151").Value = "alb000107"

$wsImages.Range("A152").Value = "pic000211"
$wsImages.Range("B152").Value = "2015-08-09-pic000211.textile"
$wsImages.Range("C152").Value = "S-15-0024.jpg"
$wsImages.Range("I152").Value = "alb000108"

# --- Albums sheet: three new album numbers ---
$wsAlbums.Range("B13").Value = "alb000106"
$wsAlbums.Range("B14").Value = "alb000107"
$wsAlbums.Range("B15").Value = "alb000108"

# --- Types and Tags sheet: normalize "textile arts" -> "textile art" ---
$wsTypes.Range("B4").Value = "textile art"

# --- View state: Objects no longer active; Images becomes active tab ---
$wsObjects.Range("H48").Select()
$wsAlbums.Range("B15").Select()
$wsTypes.Range("B4").Select()
$wsImages.Activate()
$wsImages.Range("C153").Select()
